# IFRS company_list sheet: refresh financial figures for rows 2-6 and
# drop the now-unsupported trailing data rows (7-9) down to their first
# three identifying columns only (A/B/C), matching the corrected source
# feed ("error solve ifrs list").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12) ---------------------------------------------------
$ws.Range("D2").Value  = 5775
$ws.Range("E2").Value  = 287
$ws.Range("F2").Value  = 287
$ws.Range("G2").Value  = 28
$ws.Range("H2").Value  = 36
$ws.Range("I2").Value  = -19
$ws.Range("J2").Value  = 55
$ws.Range("K2").Value  = 5005
$ws.Range("L2").Value  = 3011
$ws.Range("M2").Value  = 1994
$ws.Range("N2").Value  = 1770
$ws.Range("O2").Value  = 225
$ws.Range("P2").Value  = 149
$ws.Range("Q2").Value  = -93
$ws.Range("R2").Value  = -161
$ws.Range("S2").Value  = 220
$ws.Range("T2").Value  = 148
$ws.Range("U2").Value  = -241
$ws.Range("V2").Value  = 2028
$ws.Range("W2").Value  = 4.97
$ws.Range("X2").Value  = 0.63
$ws.Range("Y2").Value  = -1.09
$ws.Range("Z2").Value  = 0.75
$ws.Range("AA2").Value = 150.98
$ws.Range("AB2").Value = 1172.4
$ws.Range("AC2").Value = -64
$ws.Range("AD2").Value = -59.3
$ws.Range("AE2").Value = 6399
$ws.Range("AF2").Value = 0.59
$ws.Range("AG2").Value = 75
$ws.Range("AH2").Value = 1.99
$ws.Range("AI2").Value = -109.55
$ws.Range("AJ2").Value = 29742762

# --- Row 3 (2015/12) ---------------------------------------------------
$ws.Range("D3").Value  = 5995
$ws.Range("E3").Value  = 293
$ws.Range("F3").Value  = 293
$ws.Range("G3").Value  = 196
$ws.Range("H3").Value  = 132
$ws.Range("I3").Value  = 128
$ws.Range("J3").Value  = 4
$ws.Range("K3").Value  = 5285
$ws.Range("L3").Value  = 3138
$ws.Range("M3").Value  = 2147
$ws.Range("N3").Value  = 1912
$ws.Range("O3").Value  = 235
$ws.Range("P3").Value  = 149
$ws.Range("Q3").Value  = 28
$ws.Range("R3").Value  = 267
$ws.Range("S3").Value  = -205
$ws.Range("T3").Value  = 82
$ws.Range("U3").Value  = -54
$ws.Range("V3").Value  = 1903
$ws.Range("W3").Value  = 4.88
$ws.Range("X3").Value  = 2.2
$ws.Range("Y3").Value  = 6.96
$ws.Range("Z3").Value  = 2.56
$ws.Range("AA3").Value = 146.19
$ws.Range("AB3").Value = 1252.55
$ws.Range("AC3").Value = 431
$ws.Range("AD3").Value = 7.7
$ws.Range("AE3").Value = 6776
$ws.Range("AF3").Value = 0.49
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 3.02
$ws.Range("AI3").Value = 22.02
$ws.Range("AJ3").Value = 29742762

# --- Row 4 (2016/12) ---------------------------------------------------
$ws.Range("D4").Value  = 4654
$ws.Range("E4").Value  = 107
$ws.Range("F4").Value  = 107
$ws.Range("G4").Value  = 83
$ws.Range("H4").Value  = 74
$ws.Range("I4").Value  = 79
$ws.Range("J4").Value  = -5
$ws.Range("K4").Value  = 5142
$ws.Range("L4").Value  = 2933
$ws.Range("M4").Value  = 2208
$ws.Range("N4").Value  = 1972
$ws.Range("O4").Value  = 237
$ws.Range("P4").Value  = 149
$ws.Range("Q4").Value  = 84
$ws.Range("R4").Value  = 3
$ws.Range("S4").Value  = -134
$ws.Range("T4").Value  = 72
$ws.Range("U4").Value  = 12
$ws.Range("V4").Value  = 1743
$ws.Range("W4").Value  = 2.31
$ws.Range("X4").Value  = 1.58
$ws.Range("Y4").Value  = 4.07
$ws.Range("Z4").Value  = 1.41
$ws.Range("AA4").Value = 132.83
$ws.Range("AB4").Value = 1295.71
$ws.Range("AC4").Value = 266
$ws.Range("AD4").Value = 14.4
$ws.Range("AE4").Value = 6988
$ws.Range("AF4").Value = 0.55
$ws.Range("AG4").Value = 120
$ws.Range("AH4").Value = 3.13
$ws.Range("AI4").Value = 42.8
$ws.Range("AJ4").Value = 29747874

# --- Row 5 (2017/12) ---------------------------------------------------
$ws.Range("D5").Value  = 5020
$ws.Range("E5").Value  = 149
$ws.Range("F5").Value  = 149
$ws.Range("G5").Value  = 72
$ws.Range("H5").Value  = 59
$ws.Range("I5").Value  = 60
$ws.Range("J5").Value  = -1
$ws.Range("K5").Value  = 5092
$ws.Range("L5").Value  = 2663
$ws.Range("M5").Value  = 2429
$ws.Range("N5").Value  = 2210
$ws.Range("O5").Value  = 218
$ws.Range("P5").Value  = 149
$ws.Range("Q5").Value  = 281
$ws.Range("R5").Value  = 30
$ws.Range("S5").Value  = -195
$ws.Range("T5").Value  = 30
$ws.Range("U5").Value  = 251
$ws.Range("V5").Value  = 1566
$ws.Range("W5").Value  = 2.97
$ws.Range("X5").Value  = 1.18
$ws.Range("Y5").Value  = 2.88
$ws.Range("Z5").Value  = 1.16
$ws.Range("AA5").Value = 109.67
$ws.Range("AB5").Value = 1313.95
$ws.Range("AC5").Value = 202
$ws.Range("AD5").Value = 14.93
$ws.Range("AE5").Value = 7834
$ws.Range("AF5").Value = 0.39
$ws.Range("AG5").Value = 120
$ws.Range("AH5").Value = 3.97
$ws.Range("AI5").Value = 56.25
$ws.Range("AJ5").Value = 29747874

# --- Row 6 (2018/12) ----------------------------------------------------
# Note: J6/O6 were already absent before this edit (row 6 never had them).
$ws.Range("D6").Value  = 4924
$ws.Range("E6").Value  = 94
$ws.Range("F6").Value  = 94
$ws.Range("G6").Value  = -122
$ws.Range("H6").Value  = -141
$ws.Range("I6").Value  = -135
$ws.Range("K6").Value  = 4701
$ws.Range("L6").Value  = 2468
$ws.Range("M6").Value  = 2234
$ws.Range("N6").Value  = 2024
$ws.Range("P6").Value  = 149
$ws.Range("Q6").Value  = -119
$ws.Range("R6").Value  = 26
$ws.Range("S6").Value  = 68
$ws.Range("T6").Value  = 79
$ws.Range("U6").Value  = -199
$ws.Range("V6").Value  = 1699
$ws.Range("W6").Value  = 1.91
$ws.Range("X6").Value  = -2.87
$ws.Range("Y6").Value  = -6.38
$ws.Range("Z6").Value  = -2.89
$ws.Range("AA6").Value = 110.47
$ws.Range("AB6").Value = 1210.24
$ws.Range("AC6").Value = -454
$ws.Range("AD6").Value = -5.09
$ws.Range("AE6").Value = 7173
$ws.Range("AF6").Value = 0.32
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 29747874

# Dividend-yield / payout-ratio figures (AG6/AH6) are no longer reported
# for this period, so remove those two cells outright rather than zeroing.
$ws.Range("AG6:AH6").ClearContents()

# --- Rows 7-9 (2019E/2020E/2021E) ---------------------------------------
# These forward estimate rows are dropped entirely except for the leading
# identifying columns (A = index, B = "연간", C = period label).
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
